$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BRU")

# Set the PI (column C) values for rows 4, 5, 7, and 8 to match the
# new "MUSTBRUN" PI task identifiers added by the new Schuman Word
# Generator GUI.
$ws.Range("C4").Value = "MUSTBRUN2425474"
$ws.Range("C5").Value = "MUSTBRUN2425474"
$ws.Range("C7").Value = "MUSTBRUN2423960"
$ws.Range("C8").Value = "MUSTBRUN2423960"
